# Update "想去人数" (want-to-go count) figures for three events.
# Both the "展览" sheet and the "全部类型" sheet carry duplicate rows for
# these events, so both need the same updates.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 675
    $ws.Range("F8").Value = 3212
    $ws.Range("F9").Value = 4210
}
